$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 223 to make room for the newly logged "Fixing model for
# first puzzle" entry (this pushes all subsequent rows, and the references in
# the summary formula below, down by one - matching Excel's native behaviour).
$ws.Rows("223").Insert()

# Row 217: new log entry for Sarvan Amel
$ws.Range("B217").Value = "Watching Blender tutorials and creating models"
$ws.Range("C217").Value = 2.5
$ws.Range("D217").Value = "Sarvan Amel"

# Row 220: new log entry for Sarvan Amel (note trailing space preserved)
$ws.Range("B220").Value = "Creating first puzzle model "
$ws.Range("C220").Value = 3
$ws.Range("D220").Value = "Sarvan Amel"

# Row 223 (the newly inserted row): dated log entry for Sarvan Amel
$ws.Range("A223").Value = 44667
$ws.Range("B223").Value = "Fixing model for first puzzle"
$ws.Range("C223").Value = 2
$ws.Range("D223").Value = "Sarvan Amel"

# Row 226: new log entry for Sarvan Amel
$ws.Range("B226").Value = "Editing scaling of player and eviroment"
$ws.Range("C226").Value = 1
$ws.Range("D226").Value = "Sarvan Amel"

# Row 229: new log entry for Sarvan Amel
$ws.Range("B229").Value = "Fixing labirinth model"
$ws.Range("C229").Value = 2
$ws.Range("D229").Value = "Sarvan Amel"

# Row 241: Sarvan Amel's workhours total for this period
$ws.Range("B241").Value = 10.5
